$d = $word.ActiveDocument

# Re-set the text of the "Underline all edges..." paragraph to force
# recomputation / drop of the stale lastRenderedPageBreak marker.
$find = $d.Content.Find
$find.Execute("Underline all edges with the same weight as vertex key", $true, $false, $false, $false, $false, $true, 1, $false, "Underline all edges with the same weight as vertex key", 2)

Write-Output "done"
